# ---------------------------------------------------------------------------
# Add the new data row (row 4). It used to be an empty row that only carried
# an explicit (custom) row height; it now holds a real data record, using the
# same layout/styling as the existing data rows above it.
# ---------------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 254
$ws.Range("B4").Value = "PerturbRight"
$ws.Range("C4").Value = "Vehicle undergoes a perturbation to right."
$ws.Range("D4").Value = "Def/DefInvalid"

# Re-fit rows 4 and 5 so their heights go back to "automatic" (not a custom,
# pinned height) -- row 5 had no data at all, so this also drops it back to
# the sheet's default (removing the stray leftover row entry).
$ws.Rows.Item(4).AutoFit()
$ws.Rows.Item(5).AutoFit()

# ---------------------------------------------------------------------------
# Touch the very last row of the sheet so the sheet's used range / dimension
# extends all the way down to row 1048576, the same as the target file.
# ---------------------------------------------------------------------------
$ws.Range("C3").Copy()
$ws.Range("D1048576").PasteSpecial(-4122)
$ws.Range("D1048576").ClearContents()
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Move the active selection back to the top-left of the sheet, with the new
# D4 cell selected/active.
# ---------------------------------------------------------------------------
$ws.Range("D4").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
